# TrialsSetup.xlsx refresh: the "SPICE IV" trial dropped out of the
# Power Query result set on the 2025-12-22 16:00 refresh, so row 10
# (Trial Name "SPICE IV", Progress 0) is removed from the Query1
# table/worksheet and every row below it shifts up one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SPICE IV" row (worksheet row 10) - this also shrinks the
# Query1 table/autofilter range from A1:B14 to A1:B13 and shifts the
# remaining rows (ALLEGRETTO-LTE, RECOVERY, ROSETTA-Breast-01, REJOICE) up.
$ws.Rows(10).Delete()

# The hidden ExternalData_1 defined name (which records the query table's
# extent) needs to be updated to match the new, smaller range.
$wb.Names.Item("ExternalData_1").RefersTo = '=Sheet1!$A$1:$B$13'
